$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.696.99"
$ws.Range("E2").Value = "  +5.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.29"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.27"
$ws.Range("E5").Value = "  +6.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9967"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3676"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.40"
$ws.Range("E8").Value = "  +3.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3298"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.164"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("E11").Value = "  +3.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9973"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.189"
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.89"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.701.38"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.808"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001069"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06605"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.82"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9962"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.11"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.029"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("E23").Value = "  +4.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.667.90"
$ws.Range("E24").Value = "  +5.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.450"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.474"
$ws.Range("E26").Value = "  +5.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.40"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.11"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.269"
$ws.Range("E29").Value = "  +6.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.894.08"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.87"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.102"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.934"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08491"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.688"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.66"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.296"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.270"
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2117"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02255"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.448"
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6088"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.18"
$ws.Range("E44").Value = "  +12.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9964"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5842"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.62"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.998"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07215"
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.206"
$ws.Range("E51").Value = "  +3.16%  "
